$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the template-placeholder cells that are no longer used now that the
# Excel-based invoice generator has been replaced by the new HTML templates.
# (NO. PO: label in F14 is kept; only the placeholder tokens are removed.)
$ws.Range("O10").Value = ""   # was {{CUSTOMER_NAME}}
$ws.Range("O11").Value = ""   # was {{CUSTOMER_ADDRESS}}
$ws.Range("G14").Value = ""   # was {{PO_NUMBER}}

$ws.Range("F17").Value = ""   # was [[NO]]
$ws.Range("G17").Value = ""   # was [[NAME]]
$ws.Range("M17").Value = ""   # was [PCS]
$ws.Range("N17").Value = ""   # was [M2]
$ws.Range("O17").Value = ""   # was [PRICE]

$ws.Range("P26").Value = ""   # was {{SUBTOTAL}}
$ws.Range("P27").Value = ""   # was {{DP}}
$ws.Range("P28").Value = ""   # was {{ONGKIR}}
$ws.Range("P29").Value = ""   # was {{KEKURANGAN}}

# Update the saved view/selection to match the latest session.
$ws.Range("E7:Q34").Select() | Out-Null
